# Apply edits to inletConditions.xlsx as captured by the target diff.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "0.101MPa"
$ws2 = $wb.Worksheets.Item(2)   # "0.405MPa"

# --- Sheet "0.101MPa": rename the flame headers in C1:H1 to match the actual
#     Tmax values already present in row 10 (B1 "Flame1922" stays unchanged). ---
$ws1.Range("E1").Value = "Flame2063"
$ws1.Range("C1").Value = "Flame1960"
$ws1.Range("D1").Value = "Flame2013"
$ws1.Range("F1").Value = "Flame2113"
$ws1.Range("G1").Value = "Flame2160"
$ws1.Range("H1").Value = "Flame2207"

# --- Sheet "0.405MPa": fill in the missing T_in_o value for the last flame column. ---
$ws2.Range("I9").Value = 323

# --- Update sheet selections / active cells. ---
$ws1.Range("H2").Select()

# --- Make "0.405MPa" the active tab/sheet, with I13 selected. ---
$ws2.Activate()
$ws2.Range("I13").Select()
